# Scheduled runner: refresh Universalis market-price snapshot columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) for
# a batch of Leve rows across the job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 37500212
$ws.Range("I33").Value = 46154050
$ws.Range("K33").Value = 46154050
$ws.Range("M33").Value = -46153821

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 1369.9166
$ws.Range("I106").Value = 863.9
$ws.Range("K106").Value = 863.9
$ws.Range("M106").Value = -232.9

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 52490.316
$ws.Range("I113").Value = 70336.94
$ws.Range("J113").Value = 4899.3335
$ws.Range("K113").Value = 70336.94
$ws.Range("L113").Value = 4899.3335
$ws.Range("M113").Value = -67082.94
$ws.Range("N113").Value = -11407.3335

# Row 121: Mindful Medicine / Tincture of Mind
$ws.Range("H121").Value = 914
$ws.Range("J121").Value = 1078.8
$ws.Range("L121").Value = 3236.4
$ws.Range("N121").Value = -6730.4

# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1933.3334
$ws.Range("I125").Value = 1900
$ws.Range("J125").Value = 2000
$ws.Range("K125").Value = 17100
$ws.Range("L125").Value = 18000
$ws.Range("M125").Value = -14640
$ws.Range("N125").Value = -22920

# Row 131: Mindful Study / Grade 5 Tincture of Mind
$ws.Range("H131").Value = 6423.7144
$ws.Range("J131").Value = 10616.5
$ws.Range("L131").Value = 31849.5
$ws.Range("N131").Value = -41929.5

# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 10129.917
$ws.Range("I137").Value = 2893
$ws.Range("J137").Value = 12542.223
$ws.Range("K137").Value = 8679
$ws.Range("L137").Value = 37626.669
$ws.Range("M137").Value = -6129
$ws.Range("N137").Value = -42726.669

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 2980.9292
$ws.Range("I138").Value = 1234.1
$ws.Range("J138").Value = 3740.4204
$ws.Range("K138").Value = 3702.3
$ws.Range("L138").Value = 11221.2612
$ws.Range("M138").Value = 1437.7
$ws.Range("N138").Value = -21501.2612

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 4880.364
$ws.Range("I61").Value = 4409.3335
$ws.Range("K61").Value = 4409.3335
$ws.Range("M61").Value = -4197.3335

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 3272.75
$ws.Range("I132").Value = 2670.8333
$ws.Range("J132").Value = 5078.5
$ws.Range("K132").Value = 8012.499899999999
$ws.Range("L132").Value = 15235.5
$ws.Range("M132").Value = -5482.499899999999
$ws.Range("N132").Value = -20295.5

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 4880.364
$ws.Range("I136").Value = 4409.3335
$ws.Range("K136").Value = 13228.0005
$ws.Range("M136").Value = -10678.0005

$ws = $wb.Worksheets.Item("BSM")
# Row 81: Diamond Sawdust / Titanium Battleaxe
$ws.Range("H81").Value = 79780
$ws.Range("J81").Value = 79780
$ws.Range("L81").Value = 79780
$ws.Range("N81").Value = -81902

# Row 84: I'm a Lumberjack and I'm Okay (L) / Titanium Battleaxe
$ws.Range("H84").Value = 79780
$ws.Range("J84").Value = 79780
$ws.Range("L84").Value = 239340
$ws.Range("N84").Value = -249948

# Row 86: Through Thick and Thin / Adamantite Nugget
$ws.Range("H86").Value = 6723.347
$ws.Range("I86").Value = 4845.147
$ws.Range("J86").Value = 10980.6
$ws.Range("K86").Value = 4845.147
$ws.Range("L86").Value = 10980.6
$ws.Range("M86").Value = -3722.147
$ws.Range("N86").Value = -13226.6

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws.Range("H89").Value = 6723.347
$ws.Range("I89").Value = 4845.147
$ws.Range("J89").Value = 10980.6
$ws.Range("K89").Value = 24225.735
$ws.Range("L89").Value = 54903
$ws.Range("M89").Value = -18609.735
$ws.Range("N89").Value = -66135

# Row 99: Meddle in Metal / Oroshigane Ingot
$ws.Range("H99").Value = 4728.5713
$ws.Range("I99").Value = 4142.5
$ws.Range("K99").Value = 4142.5
$ws.Range("M99").Value = -2644.5

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 4657.8887
$ws.Range("I105").Value = 4740.125
$ws.Range("K105").Value = 4740.125
$ws.Range("M105").Value = -2993.125

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 3843.111
$ws.Range("J107").Value = 3750
$ws.Range("L107").Value = 3750
$ws.Range("N107").Value = -7590

# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 80132.8
$ws.Range("J138").Value = 80132.8
$ws.Range("L138").Value = 80132.8
$ws.Range("N138").Value = -90412.8

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 3639.2856
$ws.Range("I16").Value = 3496.25
$ws.Range("J16").Value = 3830
$ws.Range("K16").Value = 3496.25
$ws.Range("L16").Value = 3830
$ws.Range("M16").Value = -3209.25
$ws.Range("N16").Value = -4404

# Row 17: Say It with Spears / Feathered Harpoon
$ws.Range("H17").Value = 9
$ws.Range("J17").Value = 9
$ws.Range("L17").Value = 9
$ws.Range("N17").Value = -357

# Row 25: Bowing to Necessity / Ash Shortbow
$ws.Range("H25").Value = 805.5
$ws.Range("I25").Value = 805.5
$ws.Range("K25").Value = 805.5
$ws.Range("M25").Value = -631.5

# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 6983673
$ws.Range("I31").Value = 11949741
$ws.Range("K31").Value = 11949741
$ws.Range("M31").Value = -11949446

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 6983673
$ws.Range("I34").Value = 11949741
$ws.Range("K34").Value = 11949741
$ws.Range("M34").Value = -11949539

# Row 70: A Reward Fitting of the Faithful / Holy Cedar Necklace
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73: Just Rewards for Just Devotion (L) / Holy Cedar Necklace
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 3639.2856
$ws.Range("I113").Value = 3496.25
$ws.Range("J113").Value = 3830
$ws.Range("K113").Value = 3496.25
$ws.Range("L113").Value = 3830
$ws.Range("M113").Value = -1326.25
$ws.Range("N113").Value = -8170

# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 3510.2778
$ws.Range("I132").Value = 3261.625
$ws.Range("K132").Value = 9784.875
$ws.Range("M132").Value = -7254.875

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 13648.435
$ws.Range("I134").Value = 3057.125
$ws.Range("J134").Value = 37857.145
$ws.Range("K134").Value = 9171.375
$ws.Range("L134").Value = 113571.435
$ws.Range("M134").Value = -6636.375
$ws.Range("N134").Value = -118641.435

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 9269.286
$ws.Range("I5").Value = 593.3333
$ws.Range("J5").Value = 15776.25
$ws.Range("K5").Value = 1779.9999
$ws.Range("L5").Value = 47328.75
$ws.Range("M5").Value = -1667.9999
$ws.Range("N5").Value = -47552.75

# Row 69: Loving That Muffin Top / Ishgardian Muffin
$ws.Range("H69").Value = 580.2857
$ws.Range("I69").Value = 547.63635
$ws.Range("J69").Value = 700
$ws.Range("K69").Value = 1642.90905
$ws.Range("L69").Value = 2100
$ws.Range("M69").Value = -831.90905
$ws.Range("N69").Value = -3722

# Row 72: Muffin of the Morn (L) / Ishgardian Muffin
$ws.Range("H72").Value = 580.2857
$ws.Range("I72").Value = 547.63635
$ws.Range("J72").Value = 700
$ws.Range("K72").Value = 4928.72715
$ws.Range("L72").Value = 6300
$ws.Range("M72").Value = -872.7271499999997
$ws.Range("N72").Value = -14412

# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 9269.286
$ws.Range("I135").Value = 593.3333
$ws.Range("J135").Value = 15776.25
$ws.Range("K135").Value = 5339.9997
$ws.Range("L135").Value = 141986.25
$ws.Range("M135").Value = -2804.9997
$ws.Range("N135").Value = -147056.25

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers / Copper Ingot
$ws.Range("H2").Value = 200.84616
$ws.Range("I2").Value = 139
$ws.Range("J2").Value = 317.66666
$ws.Range("K2").Value = 139
$ws.Range("L2").Value = 317.66666
$ws.Range("M2").Value = -26
$ws.Range("N2").Value = -543.66666

# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 1600.909
$ws.Range("I107").Value = 1437.1666
$ws.Range("J107").Value = 1797.4
$ws.Range("K107").Value = 1437.1666
$ws.Range("L107").Value = 1797.4
$ws.Range("M107").Value = 482.8334
$ws.Range("N107").Value = -5637.4

# Row 109: You're My Wonderhall / Hematite Earrings of Healing
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1157
$ws.Range("I122").Value = 1157
$ws.Range("K122").Value = 3471
$ws.Range("M122").Value = -1021

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 174699.97
$ws.Range("I132").Value = 183929.08
$ws.Range("J132").Value = 5499.6665
$ws.Range("K132").Value = 551787.24
$ws.Range("L132").Value = 16498.9995
$ws.Range("M132").Value = -549257.24
$ws.Range("N132").Value = -21558.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 12065.143
$ws.Range("I7").Value = 12659.923
$ws.Range("K7").Value = 12659.923
$ws.Range("M7").Value = -12547.923

# Row 16: Saddle Sore / Hard Leather
$ws.Range("H16").Value = 8180.75
$ws.Range("I16").Value = 6648.6665
$ws.Range("J16").Value = 9100
$ws.Range("K16").Value = 6648.6665
$ws.Range("L16").Value = 9100
$ws.Range("M16").Value = -6478.6665
$ws.Range("N16").Value = -9440

# Row 98: Try Tricorne Again / Tigerskin Tricorne of Aiming
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# Row 104: Brace Yourselves / Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 47920.6
$ws.Range("J104").Value = 47920.6
$ws.Range("L104").Value = 47920.6
$ws.Range("N104").Value = -54908.6

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 12065.143
$ws.Range("I126").Value = 12659.923
$ws.Range("K126").Value = 37979.769
$ws.Range("M126").Value = -35509.769

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 4126.65
$ws.Range("I132").Value = 3469.5334
$ws.Range("J132").Value = 6098
$ws.Range("K132").Value = 10408.6002
$ws.Range("L132").Value = 18294
$ws.Range("M132").Value = -7878.600199999999
$ws.Range("N132").Value = -23354

$ws = $wb.Worksheets.Item("WVR")
# Row 92: Modest Beginnings / Bloodhempen Culottes of Casting
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 113: A Tender Table / Pixie Floss
$ws.Range("H113").Value = 777.06665
$ws.Range("I113").Value = 603
$ws.Range("K113").Value = 1809
$ws.Range("M113").Value = 361

# Row 140: Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 56748.43
$ws.Range("J140").Value = 56748.43
$ws.Range("L140").Value = 56748.43
$ws.Range("N140").Value = -67108.42999999999
